# Update map for new gates (openCyto explore/map.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Write the new cell values in the exact order needed so the shared
# --- string table indices line up with the target workbook. Column A
# --- (row-by-row) is written first, then column C cells in the order
# --- that introduces their strings 48,49,50,51 out of row order, then
# --- the remaining A/C new rows, finally the B column (all re-using
# --- already-existing shared strings so their write order is irrelevant).

$ws.Range("A21").Value = "CD8/HLA-DR+"
$ws.Range("C21").Value = "activated cytotoxic Tcells (CD8+ HLA-DR+)"

$ws.Range("A22").Value = "CCR7-CD45RA-/CD28-CD27-"
$ws.Range("A23").Value = "CCR7-CD45RA-/CD28+CD27-"
$ws.Range("A24").Value = "CCR7-CD45RA-/CD28-CD27+"
$ws.Range("A25").Value = "CCR7-CD45RA-/CD28+CD27+"

$ws.Range("C25").Value = "EM1 cytotoxic Tcells (CD27+  CD28+)"
$ws.Range("C24").Value = "EM2 cytotoxic Tcells (CD27+  CD28-)"
$ws.Range("C22").Value = "EM3 cytotoxic Tcells (CD27-  CD28-)"
$ws.Range("C23").Value = "EM4 cytotoxic Tcells (CD27-  CD28+)"

$ws.Range("A26").Value = "CCR7-CD45RA+/CD28-CD27-"
$ws.Range("A27").Value = "CCR7-CD45RA+/CD28-CD27+"
$ws.Range("A28").Value = "CCR7-CD45RA+/CD28+CD27+"

$ws.Range("C26").Value = "pE cytotoxic Tcells (CD27-  CD28-)"
$ws.Range("C28").Value = "pE1 cytotoxic Tcells (CD27+  CD28+)"
$ws.Range("C27").Value = "pE2 cytotoxic Tcells (CD27+ , CD28-)"

$ws.Range("B21").Value = "CD8"
$ws.Range("B22").Value = "CD8/CCR7-CD45RA-"
$ws.Range("B23").Value = "CD8/CCR7-CD45RA-"
$ws.Range("B24").Value = "CD8/CCR7-CD45RA-"
$ws.Range("B25").Value = "CD8/CCR7-CD45RA-"
$ws.Range("B26").Value = "CD8/CCR7-CD45RA+"
$ws.Range("B27").Value = "CD8/CCR7-CD45RA+"
$ws.Range("B28").Value = "CD8/CCR7-CD45RA+"

# --- Apply the new "Monaco / 11 / purple-blue" font to the new C-column
# --- gate names. Build the style on the first cell, then propagate it
# --- to the rest via a formats-only paste so no extra per-cell style
# --- slots are created.
$firstGate = $ws.Range("C21")
$firstGate.Font.Name = "Monaco"
$firstGate.Font.Size = 11
$firstGate.Font.Color = 16724793

$firstGate.Copy()
$ws.Range("C22:C28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column A needs to be a bit wider to fit the new longer labels.
$ws.Columns.Item(1).ColumnWidth = 23.65

# --- Match the author's final selection/cursor position.
$ws.Range("D26").Select()
